$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H101").Value = 773.06665
$ws.Range("I101").Value = 239.2
$ws.Range("J101").Value = 1040
$ws.Range("K101").Value = 717.5999999999999
$ws.Range("L101").Value = 3120
$ws.Range("M101").Value = 904.4000000000001
$ws.Range("N101").Value = -6364

$ws.Range("H112").Value = 1393.8077
$ws.Range("I112").Value = 533.3333
$ws.Range("J112").Value = 1506.0435
$ws.Range("K112").Value = 1599.9999
$ws.Range("L112").Value = 4518.1305
$ws.Range("M112").Value = -491.9999
$ws.Range("N112").Value = -6734.1305

$ws.Range("H129").Value = 888.4091
$ws.Range("J129").Value = 1124
$ws.Range("L129").Value = 3372
$ws.Range("N129").Value = -13372

$ws.Range("H134").Value = 47725
$ws.Range("J134").Value = 47725
$ws.Range("L134").Value = 47725
$ws.Range("N134").Value = -57865

$ws.Range("H136").Value = 49450
$ws.Range("J136").Value = 49450
$ws.Range("L136").Value = 49450
$ws.Range("N136").Value = -59650

$ws.Range("H137").Value = 1211.1666
$ws.Range("I137").Value = 821.23254
$ws.Range("K137").Value = 2463.69762
$ws.Range("M137").Value = 86.30238000000008

$ws.Range("H138").Value = 1858.5857
$ws.Range("I138").Value = 620.4595
$ws.Range("J138").Value = 3246.7878
$ws.Range("K138").Value = 1861.3785
$ws.Range("L138").Value = 9740.3634
$ws.Range("M138").Value = 3278.6215
$ws.Range("N138").Value = -20020.3634

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10644

$ws.Range("H61").Value = 997.8889
$ws.Range("I61").Value = 905.8387
$ws.Range("J61").Value = 1201.7142
$ws.Range("K61").Value = 905.8387
$ws.Range("L61").Value = 1201.7142
$ws.Range("M61").Value = -693.8387
$ws.Range("N61").Value = -1625.7142

$ws.Range("H74").Value = 5557814.5
$ws.Range("I74").Value = 7577752
$ws.Range("J74").Value = 2986.6667
$ws.Range("K74").Value = 7577752
$ws.Range("L74").Value = 2986.6667
$ws.Range("M74").Value = -7576878
$ws.Range("N74").Value = -4734.6667

$ws.Range("H77").Value = 5557814.5
$ws.Range("I77").Value = 7577752
$ws.Range("J77").Value = 2986.6667
$ws.Range("K77").Value = 37888760
$ws.Range("L77").Value = 14933.3335
$ws.Range("M77").Value = -37884392
$ws.Range("N77").Value = -23669.3335

$ws.Range("H136").Value = 997.8889
$ws.Range("I136").Value = 905.8387
$ws.Range("J136").Value = 1201.7142
$ws.Range("K136").Value = 2717.5161
$ws.Range("L136").Value = 3605.1426
$ws.Range("M136").Value = -167.5160999999998
$ws.Range("N136").Value = -8705.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 1393.3334
$ws.Range("I12").Value = 340
$ws.Range("J12").Value = 3500
$ws.Range("K12").Value = 340
$ws.Range("L12").Value = 3500
$ws.Range("M12").Value = -172
$ws.Range("N12").Value = -3836

$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H35").Value = 21500
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 21500
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 21500
$ws.Range("N35").Value = -22120
$ws.Range("M35").ClearContents()

$ws.Range("H86").Value = 14287664
$ws.Range("I86").Value = 22223954
$ws.Range("J86").Value = 2339.4
$ws.Range("K86").Value = 22223954
$ws.Range("L86").Value = 2339.4
$ws.Range("M86").Value = -22222831
$ws.Range("N86").Value = -4585.4

$ws.Range("H89").Value = 14287664
$ws.Range("I89").Value = 22223954
$ws.Range("J89").Value = 2339.4
$ws.Range("K89").Value = 111119770
$ws.Range("L89").Value = 11697
$ws.Range("M89").Value = -111114154
$ws.Range("N89").Value = -22929

$ws.Range("H134").Value = 1348.0817
$ws.Range("I134").Value = 1105.6279
$ws.Range("J134").Value = 3085.6667
$ws.Range("K134").Value = 3316.8837
$ws.Range("L134").Value = 9257.000100000001
$ws.Range("M134").Value = -781.8836999999999
$ws.Range("N134").Value = -14327.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 4202.2
$ws.Range("I14").Value = 3500
$ws.Range("J14").Value = 4670.3335
$ws.Range("K14").Value = 3500
$ws.Range("L14").Value = 4670.3335
$ws.Range("M14").Value = -3330
$ws.Range("N14").Value = -5010.3335

$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -713

$ws.Range("H31").Value = 3590644.5
$ws.Range("I31").Value = 5912627
$ws.Range("J31").Value = 2125.5
$ws.Range("K31").Value = 5912627
$ws.Range("L31").Value = 2125.5
$ws.Range("M31").Value = -5912332
$ws.Range("N31").Value = -2715.5

$ws.Range("H34").Value = 3590644.5
$ws.Range("I34").Value = 5912627
$ws.Range("J34").Value = 2125.5
$ws.Range("K34").Value = 5912627
$ws.Range("L34").Value = 2125.5
$ws.Range("M34").Value = -5912425
$ws.Range("N34").Value = -2529.5

$ws.Range("H36").Value = 3000
$ws.Range("I36").Value = 3000
$ws.Range("K36").Value = 3000
$ws.Range("M36").Value = -2612

$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2840

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H58").Value = 1331.0416
$ws.Range("I58").Value = 688.05554
$ws.Range("K58").Value = 688.05554
$ws.Range("M58").Value = -485.05554

$ws.Range("H132").Value = 1576.0741
$ws.Range("I132").Value = 1164.9
$ws.Range("J132").Value = 2750.8572
$ws.Range("K132").Value = 3494.7
$ws.Range("L132").Value = 8252.571599999999
$ws.Range("M132").Value = -964.7000000000003
$ws.Range("N132").Value = -13312.5716

$ws.Range("H134").Value = 2667.7812
$ws.Range("I134").Value = 3374.3684
$ws.Range("J134").Value = 1635.0769
$ws.Range("K134").Value = 10123.1052
$ws.Range("L134").Value = 4905.2307
$ws.Range("M134").Value = -7588.1052
$ws.Range("N134").Value = -9975.2307

$ws.Range("H136").Value = 1331.0416
$ws.Range("I136").Value = 688.05554
$ws.Range("K136").Value = 2064.16662
$ws.Range("M136").Value = 485.83338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 857.1429000000001
$ws.Range("I132").Value = 500
$ws.Range("J132").Value = 1333.3334
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 12000.0006
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -17060.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 3722.1428
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 4175.8335
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 4175.8335
$ws.Range("M17").Value = -832
$ws.Range("N17").Value = -4511.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 18900
$ws.Range("J21").Value = 18900
$ws.Range("L21").Value = 18900
$ws.Range("N21").Value = -19248

$ws.Range("H132").Value = 11798480
$ws.Range("I132").Value = 22330538
$ws.Range("J132").Value = 2574.08
$ws.Range("K132").Value = 66991614
$ws.Range("L132").Value = 7722.24
$ws.Range("M132").Value = -66989084
$ws.Range("N132").Value = -12782.24

$ws.Range("H136").Value = 3197.2632
$ws.Range("I136").Value = 3814.35
$ws.Range("J136").Value = 1745.2941
$ws.Range("K136").Value = 11443.05
$ws.Range("L136").Value = 5235.8823
$ws.Range("M136").Value = -8893.049999999999
$ws.Range("N136").Value = -10335.8823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 25817.5
$ws.Range("J104").Value = 25817.5
$ws.Range("L104").Value = 25817.5
$ws.Range("N104").Value = -32805.5

$ws.Range("H132").Value = 2419.8708
$ws.Range("I132").Value = 892.9474
$ws.Range("J132").Value = 4837.5
$ws.Range("K132").Value = 2678.8422
$ws.Range("L132").Value = 14512.5
$ws.Range("M132").Value = -148.8422
$ws.Range("N132").Value = -19572.5

$ws.Range("H136").Value = 1178.4546
$ws.Range("I136").Value = 623.6667
$ws.Range("J136").Value = 3675
$ws.Range("K136").Value = 1871.0001
$ws.Range("L136").Value = 11025
$ws.Range("M136").Value = 678.9999
$ws.Range("N136").Value = -16125
